$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4773
$ws.Range("I19").Value = 9397.817999999999
$ws.Range("J19").Value = 1139.2142
$ws.Range("K19").Value = 9397.817999999999
$ws.Range("L19").Value = 1139.2142
$ws.Range("M19").Value = -9222.817999999999
$ws.Range("N19").Value = -1489.2142
$ws.Range("H98").Value = 1145.2037
$ws.Range("I98").Value = 724
$ws.Range("J98").Value = 1921.1052
$ws.Range("K98").Value = 724
$ws.Range("L98").Value = 1921.1052
$ws.Range("M98").Value = 774
$ws.Range("N98").Value = -4917.1052
$ws.Range("H105").Value = 28330
$ws.Range("J105").Value = 28330
$ws.Range("L105").Value = 28330
$ws.Range("N105").Value = -35318
$ws.Range("H113").Value = 3352.2424
$ws.Range("I113").Value = 2805
$ws.Range("J113").Value = 3590.1738
$ws.Range("K113").Value = 2805
$ws.Range("L113").Value = 3590.1738
$ws.Range("M113").Value = 449
$ws.Range("N113").Value = -10098.1738
$ws.Range("H122").Value = 1145.2037
$ws.Range("I122").Value = 724
$ws.Range("J122").Value = 1921.1052
$ws.Range("K122").Value = 2172
$ws.Range("L122").Value = 5763.3156
$ws.Range("M122").Value = 278
$ws.Range("N122").Value = -10663.3156
$ws.Range("H132").Value = 2478.8135
$ws.Range("I132").Value = 1494.8541
$ws.Range("J132").Value = 6772.4546
$ws.Range("K132").Value = 4484.5623
$ws.Range("L132").Value = 20317.3638
$ws.Range("M132").Value = -1954.5623
$ws.Range("N132").Value = -25377.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 795.2
$ws.Range("I45").Value = 654
$ws.Range("J45").Value = 1007
$ws.Range("K45").Value = 654
$ws.Range("L45").Value = 1007
$ws.Range("M45").Value = -277
$ws.Range("N45").Value = -1761
$ws.Range("H61").Value = 18640.518
$ws.Range("I61").Value = 24628.488
$ws.Range("J61").Value = 3494.4707
$ws.Range("K61").Value = 24628.488
$ws.Range("L61").Value = 3494.4707
$ws.Range("M61").Value = -24416.488
$ws.Range("N61").Value = -3918.4707
$ws.Range("H74").Value = 1645.3448
$ws.Range("I74").Value = 1052.7391
$ws.Range("J74").Value = 3917
$ws.Range("K74").Value = 1052.7391
$ws.Range("L74").Value = 3917
$ws.Range("M74").Value = -178.7391
$ws.Range("N74").Value = -5665
$ws.Range("H77").Value = 1645.3448
$ws.Range("I77").Value = 1052.7391
$ws.Range("J77").Value = 3917
$ws.Range("K77").Value = 5263.6955
$ws.Range("L77").Value = 19585
$ws.Range("M77").Value = -895.6954999999998
$ws.Range("N77").Value = -28321
$ws.Range("H88").Value = 2370.7368
$ws.Range("J88").Value = 2134.2222
$ws.Range("L88").Value = 2134.2222
$ws.Range("N88").Value = -2946.2222
$ws.Range("H91").Value = 2370.7368
$ws.Range("J91").Value = 2134.2222
$ws.Range("L91").Value = 2134.2222
$ws.Range("N91").Value = -4942.2222
$ws.Range("H110").Value = 981.5
$ws.Range("I110").Value = 980.4
$ws.Range("J110").Value = 983.3333
$ws.Range("K110").Value = 980.4
$ws.Range("L110").Value = 983.3333
$ws.Range("M110").Value = 1064.6
$ws.Range("N110").Value = -5073.3333
$ws.Range("H136").Value = 18640.518
$ws.Range("I136").Value = 24628.488
$ws.Range("J136").Value = 3494.4707
$ws.Range("K136").Value = 73885.46400000001
$ws.Range("L136").Value = 10483.4121
$ws.Range("M136").Value = -71335.46400000001
$ws.Range("N136").Value = -15583.4121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4743.8945
$ws.Range("I86").Value = 7100
$ws.Range("J86").Value = 2623.4
$ws.Range("K86").Value = 7100
$ws.Range("L86").Value = 2623.4
$ws.Range("M86").Value = -5977
$ws.Range("N86").Value = -4869.4
$ws.Range("H89").Value = 4743.8945
$ws.Range("I89").Value = 7100
$ws.Range("J89").Value = 2623.4
$ws.Range("K89").Value = 35500
$ws.Range("L89").Value = 13117
$ws.Range("M89").Value = -29884
$ws.Range("N89").Value = -24349
$ws.Range("H99").Value = 4198623
$ws.Range("I99").Value = 2140581.5
$ws.Range("J99").Value = 6668273.5
$ws.Range("K99").Value = 2140581.5
$ws.Range("L99").Value = 6668273.5
$ws.Range("M99").Value = -2139083.5
$ws.Range("N99").Value = -6671269.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 803.25
$ws.Range("I16").Value = 800
$ws.Range("J16").Value = 813
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 813
$ws.Range("M16").Value = -513
$ws.Range("N16").Value = -1387
$ws.Range("H31").Value = 2634.4375
$ws.Range("I31").Value = 1998.359
$ws.Range("J31").Value = 3626.72
$ws.Range("K31").Value = 1998.359
$ws.Range("L31").Value = 3626.72
$ws.Range("M31").Value = -1703.359
$ws.Range("N31").Value = -4216.719999999999
$ws.Range("H34").Value = 2634.4375
$ws.Range("I34").Value = 1998.359
$ws.Range("J34").Value = 3626.72
$ws.Range("K34").Value = 1998.359
$ws.Range("L34").Value = 3626.72
$ws.Range("M34").Value = -1796.359
$ws.Range("N34").Value = -4030.72
$ws.Range("H58").Value = 1622.1428
$ws.Range("I58").Value = 923.8
$ws.Range("J58").Value = 2553.2666
$ws.Range("K58").Value = 923.8
$ws.Range("L58").Value = 2553.2666
$ws.Range("M58").Value = -720.8
$ws.Range("N58").Value = -2959.2666
$ws.Range("H99").Value = 69279.266
$ws.Range("I99").Value = 112852.664
$ws.Range("J99").Value = 3919.1667
$ws.Range("K99").Value = 112852.664
$ws.Range("L99").Value = 3919.1667
$ws.Range("M99").Value = -111354.664
$ws.Range("N99").Value = -6915.1667
$ws.Range("H105").Value = 1165.7693
$ws.Range("I105").Value = 708
$ws.Range("J105").Value = 1451.875
$ws.Range("K105").Value = 708
$ws.Range("L105").Value = 1451.875
$ws.Range("M105").Value = 1039
$ws.Range("N105").Value = -4945.875
$ws.Range("H107").Value = 320.6111
$ws.Range("I107").Value = 270.9565
$ws.Range("K107").Value = 270.9565
$ws.Range("M107").Value = 1649.0435
$ws.Range("H113").Value = 803.25
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 813
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 813
$ws.Range("M113").Value = 1370
$ws.Range("N113").Value = -5153
$ws.Range("H126").Value = 69279.266
$ws.Range("I126").Value = 112852.664
$ws.Range("J126").Value = 3919.1667
$ws.Range("K126").Value = 338557.992
$ws.Range("L126").Value = 11757.5001
$ws.Range("M126").Value = -336087.992
$ws.Range("N126").Value = -16697.5001
$ws.Range("H134").Value = 1488.3489
$ws.Range("I134").Value = 958.1786
$ws.Range("J134").Value = 2478
$ws.Range("K134").Value = 2874.5358
$ws.Range("L134").Value = 7434
$ws.Range("M134").Value = -339.5357999999997
$ws.Range("N134").Value = -12504
$ws.Range("H136").Value = 1622.1428
$ws.Range("I136").Value = 923.8
$ws.Range("J136").Value = 2553.2666
$ws.Range("K136").Value = 2771.4
$ws.Range("L136").Value = 7659.7998
$ws.Range("M136").Value = -221.3999999999996
$ws.Range("N136").Value = -12759.7998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1265.0526
$ws.Range("I140").Value = 855.73334
$ws.Range("J140").Value = 2800
$ws.Range("K140").Value = 2567.20002
$ws.Range("L140").Value = 8400
$ws.Range("M140").Value = 2612.79998
$ws.Range("N140").Value = -18760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1632
$ws.Range("I122").Value = 1658.0834
$ws.Range("J122").Value = 1600.7
$ws.Range("K122").Value = 4974.2502
$ws.Range("L122").Value = 4802.1
$ws.Range("M122").Value = -2524.2502
$ws.Range("N122").Value = -9702.1
$ws.Range("H126").Value = 3106.2163
$ws.Range("I126").Value = 2732.5
$ws.Range("J126").Value = 3545.8823
$ws.Range("K126").Value = 8197.5
$ws.Range("L126").Value = 10637.6469
$ws.Range("M126").Value = -5727.5
$ws.Range("N126").Value = -15577.6469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 24391702
$ws.Range("I122").Value = 33334784
$ws.Range("K122").Value = 100004352
$ws.Range("M122").Value = -100001902
$ws.Range("H136").Value = 18520610
$ws.Range("I136").Value = 24391848
$ws.Range("K136").Value = 73175544
$ws.Range("M136").Value = -73172994

Write-Host "Applied 213 cell updates"
